$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column C for rows 2 through 28
# from serial date 45552 (2024-09-17) to 45553 (2024-09-18).
for ($row = 2; $row -le 28; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45552) {
        $cell.Value2 = 45553
    }
}
